# Scheduled-runner refresh of the Leve profit calculations (currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ -> columns H:N) for
# rows whose underlying market-board snapshot changed. Values are plain numbers
# (no formulas in these sheets), so each affected cell is rewritten directly;
# a few rows also gain/lose a LeveProfit cell entirely, which is modelled as a
# ClearContents() (cell removed) vs. a fresh Value assignment (cell added).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 337.8
$ws.Range("I4").Value = 229.66667
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 229.66667
$ws.Range("L4").Value = 500
$ws.Range("M4").Value = -115.66667
$ws.Range("N4").Value = -728

$ws.Range("H33").Value = 161.22223
$ws.Range("I33").Value = 162.125
$ws.Range("J33").Value = 154
$ws.Range("K33").Value = 162.125
$ws.Range("L33").Value = 154
$ws.Range("M33").Value = 66.875

$ws.Range("H40").Value = 2492.3076
$ws.Range("I40").Value = 2116.6667
$ws.Range("J40").Value = 2814.2856
$ws.Range("K40").Value = 2116.6667
$ws.Range("L40").Value = 2814.2856
$ws.Range("M40").Value = -1941.6667
$ws.Range("N40").Value = -3164.2856

$ws.Range("H51").Value = 3799.6667
$ws.Range("I51").Value = 3799.6667
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 3799.6667
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -3315.6667
$ws.Range("N51").ClearContents()

$ws.Range("H62").Value = 7249.5
$ws.Range("I62").Value = 7249.5
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 7249.5
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -6625.5

$ws.Range("H65").Value = 7249.5
$ws.Range("I65").Value = 7249.5
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 36247.5
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -33127.5

$ws.Range("H69").Value = 11875
$ws.Range("I69").Value = 11875
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 35625
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -34751

$ws.Range("H72").Value = 11875
$ws.Range("I72").Value = 11875
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 106875
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -102507

$ws.Range("H88").Value = 4399.4
$ws.Range("I88").Value = 3996.5
$ws.Range("J88").Value = 4500.125
$ws.Range("K88").Value = 3996.5
$ws.Range("L88").Value = 4500.125
$ws.Range("M88").Value = -3590.5
$ws.Range("N88").Value = -5312.125

$ws.Range("H91").Value = 4399.4
$ws.Range("I91").Value = 3996.5
$ws.Range("J91").Value = 4500.125
$ws.Range("K91").Value = 3996.5
$ws.Range("L91").Value = 4500.125
$ws.Range("M91").Value = -2592.5
$ws.Range("N91").Value = -7308.125

$ws.Range("H107").Value = 887.8333
$ws.Range("I107").Value = 887.8333
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 887.8333
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1032.1667

$ws.Range("H132").Value = 1021.03845
$ws.Range("I132").Value = 1055.6666
$ws.Range("J132").Value = 605.5
$ws.Range("K132").Value = 3166.9998
$ws.Range("L132").Value = 1816.5
$ws.Range("M132").Value = -636.9998000000001

$ws.Range("H138").Value = 4632.557
$ws.Range("I138").Value = 3367.3872
$ws.Range("J138").Value = 5638.205
$ws.Range("K138").Value = 10102.1616
$ws.Range("L138").Value = 16914.615
$ws.Range("M138").Value = -4962.161599999999
$ws.Range("N138").Value = -27194.615

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

$ws.Range("H32").Value = 4391.827
$ws.Range("I32").Value = 2408.4443
$ws.Range("J32").Value = 17142.143
$ws.Range("K32").Value = 2408.4443
$ws.Range("L32").Value = 17142.143
$ws.Range("M32").Value = -2121.4443
$ws.Range("N32").Value = -17716.143

$ws.Range("H33").Value = 10625
$ws.Range("I33").Value = 9285.714
$ws.Range("J33").Value = 20000
$ws.Range("K33").Value = 9285.714
$ws.Range("L33").Value = 20000
$ws.Range("M33").Value = -8956.714
$ws.Range("N33").Value = -20658

$ws.Range("H61").Value = 1885.7142
$ws.Range("I61").Value = 1885.7142
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1885.7142
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1673.7142

$ws.Range("H74").Value = 1100.0769
$ws.Range("I74").Value = 1100.0769
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1100.0769
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -226.0769
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 1100.0769
$ws.Range("I77").Value = 1100.0769
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 5500.3845
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -1132.3845
$ws.Range("N77").ClearContents()

$ws.Range("H136").Value = 1885.7142
$ws.Range("I136").Value = 1885.7142
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 5657.142599999999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -3107.142599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2569.5
$ws.Range("I86").Value = 2478
$ws.Range("J86").Value = 2752.5
$ws.Range("K86").Value = 2478
$ws.Range("L86").Value = 2752.5
$ws.Range("M86").Value = -1355
$ws.Range("N86").Value = -4998.5

$ws.Range("H89").Value = 2569.5
$ws.Range("I89").Value = 2478
$ws.Range("J89").Value = 2752.5
$ws.Range("K89").Value = 12390
$ws.Range("L89").Value = 13762.5
$ws.Range("M89").Value = -6774
$ws.Range("N89").Value = -24994.5

$ws.Range("H105").Value = 3958.9
$ws.Range("I105").Value = 4456.5
$ws.Range("J105").Value = 3212.5
$ws.Range("K105").Value = 4456.5
$ws.Range("L105").Value = 3212.5
$ws.Range("M105").Value = -2709.5
$ws.Range("N105").Value = -6706.5

$ws.Range("H134").Value = 1329.85
$ws.Range("I134").Value = 1347.2106
$ws.Range("J134").Value = 1000
$ws.Range("K134").Value = 4041.6318
$ws.Range("L134").Value = 3000
$ws.Range("M134").Value = -1506.6318

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("N35").ClearContents()

$ws.Range("H41").Value = 1250
$ws.Range("I41").Value = 1250
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 1250
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -822

$ws.Range("H50").Value = 48292.332
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 48292.332
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 48292.332
$ws.Range("N50").Value = -49542.332

$ws.Range("H51").Value = 49998.668
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 49998.668
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 49998.668
$ws.Range("N51").Value = -51470.668

$ws.Range("H58").Value = 2240.1667
$ws.Range("I58").Value = 1287.05
$ws.Range("J58").Value = 7005.75
$ws.Range("K58").Value = 1287.05
$ws.Range("L58").Value = 7005.75
$ws.Range("M58").Value = -1084.05

$ws.Range("H60").Value = 11589.462
$ws.Range("I60").Value = 11589.462
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 11589.462
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -11078.462

$ws.Range("H61").Value = 49998.668
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 49998.668
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 49998.668
$ws.Range("N61").Value = -50694.668

$ws.Range("H99").Value = 14327.5
$ws.Range("I99").Value = 9739.5
$ws.Range("J99").Value = 17997.9
$ws.Range("K99").Value = 9739.5
$ws.Range("L99").Value = 17997.9
$ws.Range("M99").Value = -8241.5
$ws.Range("N99").Value = -20993.9

$ws.Range("H122").Value = 3053
$ws.Range("I122").Value = 3267.125
$ws.Range("J122").Value = 2624.75
$ws.Range("K122").Value = 9801.375
$ws.Range("L122").Value = 7874.25
$ws.Range("M122").Value = -7351.375
$ws.Range("N122").Value = -12774.25

$ws.Range("H126").Value = 14327.5
$ws.Range("I126").Value = 9739.5
$ws.Range("J126").Value = 17997.9
$ws.Range("K126").Value = 29218.5
$ws.Range("L126").Value = 53993.7
$ws.Range("M126").Value = -26748.5
$ws.Range("N126").Value = -58933.7

$ws.Range("H132").Value = 2048.4375
$ws.Range("I132").Value = 1728.9231
$ws.Range("J132").Value = 3433
$ws.Range("K132").Value = 5186.7693
$ws.Range("L132").Value = 10299
$ws.Range("M132").Value = -2656.7693

$ws.Range("H134").Value = 3557.25
$ws.Range("I134").Value = 3522.889
$ws.Range("J134").Value = 3660.3333
$ws.Range("K134").Value = 10568.667
$ws.Range("L134").Value = 10980.9999
$ws.Range("M134").Value = -8033.667000000001

$ws.Range("H136").Value = 2240.1667
$ws.Range("I136").Value = 1287.05
$ws.Range("J136").Value = 7005.75
$ws.Range("K136").Value = 3861.15
$ws.Range("L136").Value = 21017.25
$ws.Range("M136").Value = -1311.15

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").ClearContents()

$ws.Range("H80").Value = 15239.7
$ws.Range("I80").Value = 6199.6
$ws.Range("J80").Value = 24279.8
$ws.Range("K80").Value = 6199.6
$ws.Range("L80").Value = 24279.8
$ws.Range("M80").Value = -5201.6
$ws.Range("N80").Value = -26275.8

$ws.Range("H83").Value = 15239.7
$ws.Range("I83").Value = 6199.6
$ws.Range("J83").Value = 24279.8
$ws.Range("K83").Value = 30998
$ws.Range("L83").Value = 121399
$ws.Range("M83").Value = -26006
$ws.Range("N83").Value = -131383

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2678.4167
$ws.Range("I22").Value = 1934.25
$ws.Range("J22").Value = 4166.75
$ws.Range("K22").Value = 1934.25
$ws.Range("L22").Value = 4166.75
$ws.Range("M22").Value = -1639.25

$ws.Range("H27").Value = 2678.4167
$ws.Range("I27").Value = 1934.25
$ws.Range("J27").Value = 4166.75
$ws.Range("K27").Value = 1934.25
$ws.Range("L27").Value = 4166.75
$ws.Range("M27").Value = -1827.25

$ws.Range("H68").Value = 4650.75
$ws.Range("I68").Value = 3925
$ws.Range("J68").Value = 5376.5
$ws.Range("K68").Value = 3925
$ws.Range("L68").Value = 5376.5
$ws.Range("M68").Value = -3176
$ws.Range("N68").Value = -6874.5

$ws.Range("H71").Value = 4650.75
$ws.Range("I71").Value = 3925
$ws.Range("J71").Value = 5376.5
$ws.Range("K71").Value = 19625
$ws.Range("L71").Value = 26882.5
$ws.Range("M71").Value = -15881
$ws.Range("N71").Value = -34370.5

$ws.Range("H136").Value = 4498
$ws.Range("I136").Value = 4498
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 13494
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -10944

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 1000
$ws.Range("I43").Value = 1000
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 1000
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -851

$ws.Range("H82").Value = 40000
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 40000
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 40000
$ws.Range("N82").Value = -40766

$ws.Range("H85").Value = 40000
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 40000
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 40000
$ws.Range("N85").Value = -42652

$ws.Range("H122").Value = 6000
$ws.Range("I122").Value = 6000
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 18000
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -15550

$ws.Range("H136").Value = 1701.6666
$ws.Range("I136").Value = 1751.8572
$ws.Range("J136").Value = 999
$ws.Range("K136").Value = 5255.571599999999
$ws.Range("L136").Value = 2997
$ws.Range("M136").Value = -2705.571599999999
